# Scheduled market-data refresh: update current average prices / profit
# calcs for the leves that received fresh data this run, and clear the
# stale H:N (price/profit) columns for leves that aged out without a
# fresh quote.

$wb = $excel.ActiveWorkbook

# ---- ALC: row 131 gets a refreshed price quote -----------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3499.75
$ws.Range("I131").Value = 999
$ws.Range("K131").Value = 2997
$ws.Range("M131").Value = 2043

# ---- CRP: row 4 price data expired to zero; rows 129-141 (minus 136) -
# ---- lost their stale quotes entirely ---------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$crpClearRows = @(129, 130, 131, 132, 133, 134, 135, 137, 138, 139, 140, 141)
foreach ($r in $crpClearRows) {
    $ws.Range("H" + $r + ":N" + $r).ClearContents()
}

# ---- LTW: rows 7 & 126 get refreshed quotes; rows 124-141 (minus 126) -
# ---- lost their stale quotes entirely ---------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 817
$ws.Range("I7").Value = 779.4
$ws.Range("K7").Value = 779.4
$ws.Range("M7").Value = -667.4

$ws.Range("H126").Value = 817
$ws.Range("I126").Value = 779.4
$ws.Range("K126").Value = 2338.2
$ws.Range("M126").Value = 131.8000000000002

$ltwClearRows = @(124, 125, 127, 128, 129, 130, 131, 132, 133, 134, 135, 136, 137, 138, 139, 140, 141)
foreach ($r in $ltwClearRows) {
    $ws.Range("H" + $r + ":N" + $r).ClearContents()
}

# ---- WVR: rows 62 & 65 get refreshed quotes ----------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6009.273
$ws.Range("I62").Value = 2450.6667
$ws.Range("J62").Value = 7343.75
$ws.Range("K62").Value = 2450.6667
$ws.Range("L62").Value = 7343.75
$ws.Range("M62").Value = -1826.6667
$ws.Range("N62").Value = -8591.75

$ws.Range("H65").Value = 6009.273
$ws.Range("I65").Value = 2450.6667
$ws.Range("J65").Value = 7343.75
$ws.Range("K65").Value = 12253.3335
$ws.Range("L65").Value = 36718.75
$ws.Range("M65").Value = -9133.333500000001
$ws.Range("N65").Value = -42958.75
